$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.589.68'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '1.665.51'
$ws.Range('E3').Value = '  -3.54%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.514'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.16%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.62'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.262'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0882'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('D12').Value = '1.902.23'
$ws.Range('E12').Value = '  -3.47%  '
$ws.Range('D13').Value = '1.656.69'
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.15'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.54%  '
$ws.Range('E15').Value = '  -1.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '247.05'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('D18').Value = '27.608.19'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('E19').Value = '  -3.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.63%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E22').Value = '  -3.64%  '
$ws.Range('E23').Value = '  -5.07%  '
$ws.Range('E24').Value = '  -4.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.08'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.39'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.59%  '
$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.112'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.54%  '
$ws.Range('E30').Value = '  +3.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0506'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.11%  '
$ws.Range('E32').Value = '  -2.97%  '
$ws.Range('D33').Value = '1.479.51'
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('E34').Value = '  -5.20%  '
$ws.Range('E35').Value = '  -6.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.937'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.78%  '
$ws.Range('E37').Value = '  -1.10%  '
$ws.Range('E38').Value = '  -6.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0171'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.63'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('E41').Value = '  -5.63%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  -7.61%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.809.11'
$ws.Range('E44').Value = '  -3.47%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.788'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('E47').Value = '  -3.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.24'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('D49').Value = '0.0₆0108'
$ws.Range('E49').Value = '  -2.63%  '
$ws.Range('E50').Value = '  -3.05%  '
$ws.Range('E51').Value = '  -4.58%  '
